# Apply the PR #223 "Pytest" changes to Config.xlsx:
#  - Constants sheet: add two new rows
#       * MaxExecutionAttemptsHigh / 99999 / "Maximum number of execution
#         attempts for a process step which by default is high."
#       * RetryIntervalLow / 1 / "Duration, in seconds, between
#         re-execution attempts" (same description text used by RetryInterval)
#  - Settings sheet: move the active-cell selection from A27 to A24
#  - Constants sheet: move the active-cell selection to the newly
#    inserted row (entire row 10 selected, active cell A10)

$wb = $excel.ActiveWorkbook

$settings  = $wb.Worksheets.Item("Settings")
$constants = $wb.Worksheets.Item("Constants")

# --- Settings sheet: just update the remembered selection -----------------
$settings.Activate()
$settings.Range("A24").Select()

# --- Constants sheet: insert the two new constant rows --------------------
$constants.Activate()

# New row 7: MaxExecutionAttemptsHigh (inserted above the existing
# "MaxLockTimeout" row, pushing it and everything below down by one)
$constants.Rows.Item(7).Insert()
$constants.Rows.Item(7).RowHeight = 14.25
$constants.Range("A7").Value() = "MaxExecutionAttemptsHigh"
$constants.Range("B7").Value() = 99999
$constants.Range("C7").Value() = "Maximum number of execution attempts for a process step which by default is high."

# New row 10: RetryIntervalLow (inserted right after "RetryInterval", which
# is now on row 9, and before the blank separator row)
$constants.Rows.Item(10).Insert()
$constants.Rows.Item(10).RowHeight = 14.25
$constants.Range("A10").Value() = "RetryIntervalLow"
$constants.Range("B10").Value() = 1
$constants.Range("C10").Value() = "Duration, in seconds, between re-execution attempts"

# Leave the worksheet selection on the newly-added row, matching the
# recorded UI state (whole row 10 selected, active cell A10)
$constants.Rows.Item(10).Select()
